$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 304.1111
$ws.Range("I12").Value = 197.85715
$ws.Range("K12").Value = 197.85715
$ws.Range("M12").Value = -27.85714999999999

$ws.Range("H70").Value = 4259.9
$ws.Range("I70").Value = 3300
$ws.Range("J70").Value = 4499.875
$ws.Range("K70").Value = 9900
$ws.Range("L70").Value = 13499.625
$ws.Range("M70").Value = -9630
$ws.Range("N70").Value = -14039.625

$ws.Range("H73").Value = 4259.9
$ws.Range("I73").Value = 3300
$ws.Range("J73").Value = 4499.875
$ws.Range("K73").Value = 9900
$ws.Range("L73").Value = 13499.625
$ws.Range("M73").Value = -8964
$ws.Range("N73").Value = -15371.625

$ws.Range("H86").Value = 11165.777
$ws.Range("I86").Value = 4918.4
$ws.Range("J86").Value = 18975
$ws.Range("K86").Value = 4918.4
$ws.Range("L86").Value = 18975
$ws.Range("M86").Value = -3795.4
$ws.Range("N86").Value = -21221

$ws.Range("H89").Value = 11165.777
$ws.Range("I89").Value = 4918.4
$ws.Range("J89").Value = 18975
$ws.Range("K89").Value = 24592
$ws.Range("L89").Value = 94875
$ws.Range("M89").Value = -18976
$ws.Range("N89").Value = -106107

$ws.Range("H106").Value = 7799.5
$ws.Range("I106").Value = 5599
$ws.Range("J106").Value = 10000
$ws.Range("K106").Value = 5599
$ws.Range("L106").Value = 10000
$ws.Range("M106").Value = -4968
$ws.Range("N106").Value = -11262

$ws.Range("H138").Value = 2352.7144
$ws.Range("I138").Value = 882
$ws.Range("K138").Value = 2646
$ws.Range("M138").Value = 2494

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 6961.5
$ws.Range("I61").Value = 6961.5
$ws.Range("K61").Value = 6961.5
$ws.Range("M61").Value = -6749.5

$ws.Range("H63").Value = 2440.5454
$ws.Range("I63").Value = 1793.25
$ws.Range("J63").Value = 4166.6665
$ws.Range("K63").Value = 1793.25
$ws.Range("L63").Value = 4166.6665
$ws.Range("M63").Value = -1107.25
$ws.Range("N63").Value = -5538.6665

$ws.Range("H66").Value = 2440.5454
$ws.Range("I66").Value = 1793.25
$ws.Range("J66").Value = 4166.6665
$ws.Range("K66").Value = 8966.25
$ws.Range("L66").Value = 20833.3325
$ws.Range("M66").Value = -5534.25
$ws.Range("N66").Value = -27697.3325

$ws.Range("H132").Value = 1885.4
$ws.Range("J132").Value = 2257
$ws.Range("L132").Value = 6771
$ws.Range("N132").Value = -11831

$ws.Range("H136").Value = 6961.5
$ws.Range("I136").Value = 6961.5
$ws.Range("K136").Value = 20884.5
$ws.Range("M136").Value = -18334.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1000
$ws.Range("I20").Value = 1000
$ws.Range("K20").Value = 1000
$ws.Range("M20").Value = -753

$ws.Range("H86").Value = 2075.2942
$ws.Range("I86").Value = 1523.1666
$ws.Range("J86").Value = 2376.4546
$ws.Range("K86").Value = 1523.1666
$ws.Range("L86").Value = 2376.4546
$ws.Range("M86").Value = -400.1666
$ws.Range("N86").Value = -4622.4546

$ws.Range("H89").Value = 2075.2942
$ws.Range("I89").Value = 1523.1666
$ws.Range("J89").Value = 2376.4546
$ws.Range("K89").Value = 7615.833000000001
$ws.Range("L89").Value = 11882.273
$ws.Range("M89").Value = -1999.833000000001
$ws.Range("N89").Value = -23114.273

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4852.826
$ws.Range("I31").Value = 3915.5
$ws.Range("J31").Value = 6310.8887
$ws.Range("K31").Value = 3915.5
$ws.Range("L31").Value = 6310.8887
$ws.Range("M31").Value = -3620.5
$ws.Range("N31").Value = -6900.8887

$ws.Range("H34").Value = 4852.826
$ws.Range("I34").Value = 3915.5
$ws.Range("J34").Value = 6310.8887
$ws.Range("K34").Value = 3915.5
$ws.Range("L34").Value = 6310.8887
$ws.Range("M34").Value = -3713.5
$ws.Range("N34").Value = -6714.8887

$ws.Range("H39").Value = 50244.75
$ws.Range("J39").Value = 49993
$ws.Range("L39").Value = 49993
$ws.Range("N39").Value = -50775

$ws.Range("H49").Value = 50244.75
$ws.Range("J49").Value = 49993
$ws.Range("L49").Value = 49993
$ws.Range("N49").Value = -50357

$ws.Range("H99").Value = 5762.5
$ws.Range("I99").Value = 8900
$ws.Range("J99").Value = 2625
$ws.Range("K99").Value = 8900
$ws.Range("L99").Value = 2625
$ws.Range("M99").Value = -7402
$ws.Range("N99").Value = -5621

$ws.Range("H126").Value = 5762.5
$ws.Range("I126").Value = 8900
$ws.Range("J126").Value = 2625
$ws.Range("K126").Value = 26700
$ws.Range("L126").Value = 7875
$ws.Range("M126").Value = -24230
$ws.Range("N126").Value = -12815

$ws.Range("H132").Value = 1333.3684
$ws.Range("I132").Value = 1038.2
$ws.Range("K132").Value = 3114.6
$ws.Range("M132").Value = -584.6000000000004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 97500
$ws.Range("J37").Value = 97500
$ws.Range("L37").Value = 292500
$ws.Range("N37").Value = -292724

$ws.Range("H103").Value = 308.2
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

$ws.Range("H136").Value = 0
$ws.Range("J136").Value = 0
$ws.Range("L136").Value = 0
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 6946395.5
$ws.Range("I122").Value = 7354680
$ws.Range("K122").Value = 22064040
$ws.Range("M122").Value = -22061590

$ws.Range("H132").Value = 2072.8572
$ws.Range("I132").Value = 2072.8572
$ws.Range("K132").Value = 6218.571599999999
$ws.Range("M132").Value = -3688.571599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2769.9092
$ws.Range("I22").Value = 2785.7144
$ws.Range("K22").Value = 2785.7144
$ws.Range("M22").Value = -2490.7144

$ws.Range("H27").Value = 2769.9092
$ws.Range("I27").Value = 2785.7144
$ws.Range("K27").Value = 2785.7144
$ws.Range("M27").Value = -2678.7144

$ws.Range("H55").Value = 400
$ws.Range("I55").Value = 333.33334
$ws.Range("J55").Value = 466.66666
$ws.Range("K55").Value = 333.33334
$ws.Range("L55").Value = 466.66666
$ws.Range("M55").Value = -160.33334
$ws.Range("N55").Value = -812.66666

$ws.Range("H132").Value = 7517.3
$ws.Range("I132").Value = 5234.6
$ws.Range("K132").Value = 15703.8
$ws.Range("M132").Value = -13173.8

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 18265.166
$ws.Range("I113").Value = 26398.5
$ws.Range("K113").Value = 79195.5
$ws.Range("M113").Value = -77025.5

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

$ws.Range("H136").Value = 1750
$ws.Range("I136").Value = 1500
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 4500
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -1950
$ws.Range("N136").Value = -11100
